# Apply updated predicted prices and compared returns values
# to columns G (Return_with_prediction), H (return_pct_change), and I (mean_return_pct_change, row 2 only)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1233454777090742
$ws.Range("H2").Value = -6.144588705300275
$ws.Range("I2").Value = 10.03916008152928
$ws.Range("G3").Value = 0.1164891597298832
$ws.Range("H3").Value = 30.83201692794194
$ws.Range("G4").Value = -0.5882086705124271
$ws.Range("H4").Value = 4.409694460364111
$ws.Range("G5").Value = -0.6085178130136797
$ws.Range("H5").Value = 0.3399041458321979
$ws.Range("G6").Value = 0.1725906262255135
$ws.Range("H6").Value = -29.87193507018294
$ws.Range("G7").Value = 0.3937975325801231
$ws.Range("H7").Value = 140.3814324093497
$ws.Range("G8").Value = 0.1593147625250467
$ws.Range("H8").Value = -3.626038785860626
$ws.Range("G9").Value = 0.230044738459763
$ws.Range("H9").Value = 17.91570285563228
$ws.Range("G10").Value = -0.129448245585633
$ws.Range("H10").Value = -126.5205930880392
$ws.Range("G11").Value = -0.1361811319077837
$ws.Range("H11").Value = -14.66115580843429
$ws.Range("G12").Value = 0.1651443927579213
$ws.Range("H12").Value = 3.845270972188648
$ws.Range("G13").Value = 0.2732445103411564
$ws.Range("H13").Value = 32.86112460306565
$ws.Range("G14").Value = 0.1670501791011007
$ws.Range("H14").Value = -11.78778823504936
$ws.Range("G15").Value = 0.21491546478028
$ws.Range("H15").Value = -14.0012782322312
$ws.Range("G16").Value = -0.01353604261625696
$ws.Range("H16").Value = -137.1051288533722
$ws.Range("G17").Value = 0.005065242007242547
$ws.Range("H17").Value = -85.71988230133326
$ws.Range("G18").Value = 0.09163839213934212
$ws.Range("H18").Value = -47.12586748308052
$ws.Range("G19").Value = 0.0890493099100888
$ws.Range("H19").Value = -29.18709719168035
$ws.Range("G20").Value = 0.04115238075865552
$ws.Range("H20").Value = -64.10670509954818
$ws.Range("G21").Value = 0.08526401494366023
$ws.Range("H21").Value = -15.07468922642086
$ws.Range("G22").Value = 0.06602384018274093
$ws.Range("H22").Value = -29.90886157730819
$ws.Range("G23").Value = 0.1014958325285188
$ws.Range("H23").Value = -6.446821016638834
$ws.Range("G24").Value = -0.2466159566609749
$ws.Range("H24").Value = -97.88746086205579
$ws.Range("G25").Value = -0.2651334617833642
$ws.Range("H25").Value = -19.18570060398115
$ws.Range("G26").Value = 0.1678981730777884
$ws.Range("H26").Value = 5.608412015633782
$ws.Range("G27").Value = 0.1525105849795829
$ws.Range("H27").Value = -23.92209699251315
$ws.Range("G28").Value = 0.03625025585004506
$ws.Range("H28").Value = 551.4100069019112
$ws.Range("G29").Value = 0.05470951859184871
$ws.Range("H29").Value = 255.7766061239308
